$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.474.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.32%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.804.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.93%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'227.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +2.78%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.18%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'36.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.54%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.48%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -0.44%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +1.19%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'2.064.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.93%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.85%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.812.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.46%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.646"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.15%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'4.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.70%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'34.432.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.46%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'70.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.60%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'245.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.87%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -1.48%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +0.42%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.28%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.05%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +3.37%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'172.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.29%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'8.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +8.58%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'16.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.92%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +1.12%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +0.12%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.20%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.20%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.41%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.51%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -2.29%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.392.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.69%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.674"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.77%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -6.48%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -0.36%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.73%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'82.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.00%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.58%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -0.84%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +0.66%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +7.89%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.80%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'6.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.17%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -4.37%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.965.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.99%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'104.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.47%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +0.25%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -2.52%  "
$ws.Range("E51").Style = "Normal"
